$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: new/changed field names (option_02, option_03 enter the shared-string
# table first, matching the commit's string order)
$ws.Range("E2").Value = "option_02"
$ws.Range("F2").Value = "option_03"

# Row 1: new header cell G1 = "n" (type marker for the new numeric column)
$ws.Range("G1").Value = "n"

# Row 2 (continued): wrongOption field name
$ws.Range("G2").Value = "wrongOption"

# Row 3: G3 becomes numeric 0, old text ("//sound도 index랑 동일하게 ") moves to H3
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = "//sound도 index랑 동일하게 "

# Row 4: new numeric cell G4 = 1
$ws.Range("G4").Value = 1

# Row 5: new numeric cell G5 = 2
$ws.Range("G5").Value = 2

# Update selection to match target (K8)
$ws.Range("K8").Select()
